$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.9277329896308137, 1.1445802505529814, 0.48324801229439407, 2.57378440550603)"
$ws.Range("C2").Value = "NIG(1.6739832197994373, 1.2822168749098057, 3.1515433113175737, 6.2619643830078235)"
$ws.Range("D2").Value = "NIG(1.0245737617126567, 0.6081987742658426, 2.0892113614187715, 3.540535309185457)"
$ws.Range("E2").Value = "EXN(1.8865462401700384, 3.413986770425657, 3.7602501973453553)"
